$d = $word.ActiveDocument

# Step 1: Replace "de référence à une variable" with "conditionnelles" so
# the sentence reads "...les balises conditionnelles : ". This is a plain
# text substitution; at this point the paragraph is still a single run.
$d.Content.Find.Execute("de référence à une variable", $true, $false, $false, $false, $false, $true, 1, $false, "conditionnelles", 2)

# Step 2: Locate the word "conditionnelles" inside the (now updated)
# first paragraph and force it into its own run, distinct from the text
# before and after it, by nudging a direct-character-formatting property
# on exactly that sub-range and then putting it back - this causes the
# run to be split at both boundaries without altering the visible
# formatting (Bold ends up False again, matching the surrounding runs).
$p = $d.Paragraphs.First
$word1 = $p.Range.Duplicate
$word1.Find.Execute("conditionnelles")

$word1.Font.Bold = $true
$word1.Font.Bold = $false
